{"js": "const body = context.document.body;\nconst replacements = [\n  [\"2024-04-07 Sunday\", \"2024-04-08 Monday\"],\n  [\"79\u00d757=4503\", \"36\u00d730=1080\"],\n  [\"39\u00d747=1833\", \"16\u00d743=688\"],\n  [\"97\u00d782=7954\", \"51\u00d769=3519\"],\n  [\"76\u00d719=1444\", \"29\u00d787=2523\"],\n  [\"17\u00d716=272\", \"44\u00d767=2948\"],\n  [\"97\u00d723=2231\", \"56\u00d773=4088\"],\n  [\"82\u00d775=6150\", \"88\u00d724=2112\"],\n  [\"98\u00d742=4116\", \"12\u00d786=1032\"],\n  [\"21\u00d743=903\", \"33\u00d732=1056\"],\n  [\"14\u00d738=532\", \"61\u00d753=3233\"],\n  [\"12\u00d739=468\", \"71\u00d764=4544\"],\n  [\"81\u00d760=4860\", \"73\u00d713=949\"],\n  [\"39\u00d725=975\", \"67\u00d770=4690\"],\n  [\"87\u00d784=7308\", \"71\u00d786=6106\"],\n  [\"19\u00d795=1805\", \"58\u00d776=4408\"],\n  [\"32\u00d722=704\", \"26\u00d772=1872\"],\n  [\"67\u00d799=6633\", \"88\u00d737=3256\"],\n  [\"73\u00d753=3869\", \"91\u00d745=4095\"],\n  [\"21\u00d729=609\", \"16\u00d745=720\"],\n  [\"98\u00d790=8820\", \"64\u00d762=3968\"],\n  [\"52\u00d713=676\", \"21\u00d769=1449\"],\n  [\"45\u00d790=4050\", \"60\u00d722=1320\"],\n  [\"79\u00d770=5530\", \"55\u00d719=1045\"],\n  [\"31\u00d723=713\", \"23\u00d776=1748\"],\n  [\"26\u00d783=2158\", \"81\u00d725=2025\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$replacements = @(\n    @{ Old = \"2024-04-07 Sunday\"; New = \"2024-04-08 Monday\" }\n    @{ Old = \"79\u00d757=4503\"; New = \"36\u00d730=1080\" }\n    @{ Old = \"39\u00d747=1833\"; New = \"16\u00d743=688\" }\n    @{ Old = \"97\u00d782=7954\"; New = \"51\u00d769=3519\" }\n    @{ Old = \"76\u00d719=1444\"; New = \"29\u00d787=2523\" }\n    @{ Old = \"17\u00d716=272\"; New = \"44\u00d767=2948\" }\n    @{ Old = \"97\u00d723=2231\"; New = \"56\u00d773=4088\" }\n    @{ Old = \"82\u00d775=6150\"; New = \"88\u00d724=2112\" }\n    @{ Old = \"98\u00d742=4116\"; New = \"12\u00d786=1032\" }\n    @{ Old = \"21\u00d743=903\"; New = \"33\u00d732=1056\" }\n    @{ Old = \"14\u00d738=532\"; New = \"61\u00d753=3233\" }\n    @{ Old = \"12\u00d739=468\"; New = \"71\u00d764=4544\" }\n    @{ Old = \"81\u00d760=4860\"; New = \"73\u00d713=949\" }\n    @{ Old = \"39\u00d725=975\"; New = \"67\u00d770=4690\" }\n    @{ Old = \"87\u00d784=7308\"; New = \"71\u00d786=6106\" }\n    @{ Old = \"19\u00d795=1805\"; New = \"58\u00d776=4408\" }\n    @{ Old = \"32\u00d722=704\"; New = \"26\u00d772=1872\" }\n    @{ Old = \"67\u00d799=6633\"; New = \"88\u00d737=3256\" }\n    @{ Old = \"73\u00d753=3869\"; New = \"91\u00d745=4095\" }\n    @{ Old = \"21\u00d729=609\"; New = \"16\u00d745=720\" }\n    @{ Old = \"98\u00d790=8820\"; New = \"64\u00d762=3968\" }\n    @{ Old = \"52\u00d713=676\"; New = \"21\u00d769=1449\" }\n    @{ Old = \"45\u00d790=4050\"; New = \"60\u00d722=1320\" }\n    @{ Old = \"79\u00d770=5530\"; New = \"55\u00d719=1045\" }\n    @{ Old = \"31\u00d723=713\"; New = \"23\u00d776=1748\" }\n    @{ Old = \"26\u00d783=2158\"; New = \"81\u00d725=2025\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll)\n}"}
